# Add a new queried word ("NOTHING") to the "Words to find" sheet, in its
# correct alphabetically-sorted position (row 12, between CHERISH and
# DECEMBER), pushing the following entries down by one row. Also leave the
# "Words to find" sheet selected/active (as the last thing the author did
# before saving), matching the tab/selection state recorded in the file.

$wb  = $excel.ActiveWorkbook
$wsWords = $wb.Worksheets.Item("Words to find")

# Insert a new row above row 12 and fill it with the missing word so the
# list stays alphabetically sorted (CHERISH, NOTHING, DECEMBER, ...).
$wsWords.Range("A12").EntireRow.Insert()
$wsWords.Range("A12").Value = "NOTHING"

# Make "Words to find" the active sheet/tab and select A12, matching the
# view state captured in the saved workbook.
$wsWords.Activate()
$wsWords.Range("A12").Select()
